# feat: add 2022-Q3 data
#
# The workbook has a "总计" (totals) sheet and a single "2022-Q2"
# fund-holdings sheet. This adds a new "2022-Q3" fund-holdings sheet
# (placed between "总计" and "2022-Q2", newest quarter first) and a
# matching summary row on "总计", while leaving the original "2022-Q2"
# sheet's data untouched.
#
# To land on the same sheetId/rId numbering as the authored change
# (总计=1, 2022-Q3=2, 2022-Q2=3) the existing "2022-Q2" worksheet is
# renamed in place to "2022-Q3" and repopulated with the Q3 figures
# (so it keeps sheetId 2), and a fresh worksheet is appended right
# after it, named "2022-Q2", and repopulated with the original Q2
# figures (so it mints the new, higher sheetId 3).

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Write a literal text value that must not be coerced to a number/date
    # (fund codes like "004685", ratios like "15.28") while leaving the
    # cell's style index as plain "Normal" afterwards.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

function Write-FundSheet($sheet, $headerStyleSrc, $aColStyleSrc, $rows) {
    $sheet.Range("B1").Value = "基金代码"
    $sheet.Range("C1").Value = "基金名称"
    $sheet.Range("D1").Value = "基金规模"
    $sheet.Range("E1").Value = "股票总仓位"
    $sheet.Range("F1").Value = "仓位占比"
    $sheet.Range("G1").Value = "持有市值(亿元)"
    $sheet.Range("H1").Value = "仓位排名"

    $r = 2
    foreach ($row in $rows) {
        $sheet.Range("A$r").Value = $r - 2
        Set-TextValue $sheet.Range("B$r") $row[0]
        Set-TextValue $sheet.Range("C$r") $row[1]
        Set-TextValue $sheet.Range("D$r") $row[2]
        Set-TextValue $sheet.Range("E$r") $row[3]
        Set-TextValue $sheet.Range("F$r") $row[4]
        Set-TextValue $sheet.Range("G$r") $row[5]
        $sheet.Range("H$r").Value = $row[6]
        $r++
    }
    $lastRow = $r - 1

    $headerStyleSrc.Copy()
    $sheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null
    $aColStyleSrc.Copy()
    $sheet.Range("A2:A$lastRow").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

$total = $wb.Worksheets.Item("总计")
$oldQ2 = $wb.Worksheets.Item("2022-Q2")

# Original 2022-Q2 fund-holdings rows, read off the sheet before it gets
# renamed/overwritten below: [code, name, scale, stockPosTotal, posRatio,
# marketValue, posRank].
$q2Rows = @(
    , @("005585", "银河文体娱乐主题灵活配置混合", "4.98", "92.64", "4.25", "0.2116", 10)
    , @("004685", "金元顺安元启灵活配置混合", "12.44", "75.65", "0.78", "0.0970", 10)
)
$q2HeaderStyleSrc = $oldQ2.Range("B1")
$q2AColStyleSrc = $oldQ2.Range("A2")

# New 2022-Q3 fund-holdings row.
$q3Rows = @(
    , @("004685", "金元顺安元启灵活配置混合", "15.28", "77.14", "0.86", "0.1314", 5)
)

# --- Step 1: the live "2022-Q2" sheet becomes "2022-Q3" --------------------
$q3 = $oldQ2
$q3.Cells.Clear()
$q3.Name = "2022-Q3"
# The new Q3 sheet takes its header/A-column formatting from the "总计"
# sheet (both end up on cellXfs style index 2 in the authored workbook).
Write-FundSheet $q3 $total.Range("B1") $total.Range("A2") $q3Rows

# --- Step 2: a brand-new sheet becomes "2022-Q2", with the original data --
$q2 = $wb.Worksheets.Add($null, $q3)
$q2.Name = "2022-Q2"
Write-FundSheet $q2 $q2HeaderStyleSrc $q2AColStyleSrc $q2Rows

# --- Step 3: update the "总计" summary sheet --------------------------------
# Existing Q2 summary row moves down to row 3; new Q3 summary row takes
# row 2 (newest quarter first).
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.31

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.13

$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
